# Scheduled-runner update: refresh cached market-price / profit figures
# (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit
# sheets. Values below were recomputed upstream; we just push the new
# cached numbers into the same cells (a few rows' H/J/L revert to 0 with
# their dependent N/M cells cleared, and one row gains a new N cell).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2750
$ws.Range("J43").Value = 2750
$ws.Range("L43").Value = 2750
$ws.Range("N43").Value = -2888

$ws.Range("H58").Value = 973.5
$ws.Range("I58").Value = 526.1111
$ws.Range("K58").Value = 1578.3333
$ws.Range("M58").Value = -1428.3333

$ws.Range("H61").Value = 1111.1666
$ws.Range("I61").Value = 1111.1666
$ws.Range("K61").Value = 3333.4998
$ws.Range("M61").Value = -3161.4998

$ws.Range("H74").Value = 7099.815
$ws.Range("I74").Value = 3828.7
$ws.Range("K74").Value = 3828.7
$ws.Range("M74").Value = -2892.7

$ws.Range("H77").Value = 7099.815
$ws.Range("I77").Value = 3828.7
$ws.Range("K77").Value = 19143.5
$ws.Range("M77").Value = -14463.5

$ws.Range("H114").Value = 104828.5
$ws.Range("J114").Value = 104828.5
$ws.Range("L114").Value = 104828.5
$ws.Range("N114").Value = -113506.5

$ws.Range("H132").Value = 78291
$ws.Range("I132").Value = 88320.25999999999
$ws.Range("K132").Value = 264960.78
$ws.Range("M132").Value = -262430.78

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5252.676
$ws.Range("I2").Value = 970.2222
$ws.Range("K2").Value = 970.2222
$ws.Range("M2").Value = -857.2222

$ws.Range("H61").Value = 878842.5
$ws.Range("I61").Value = 927417.1
$ws.Range("K61").Value = 927417.1
$ws.Range("M61").Value = -927205.1

$ws.Range("H74").Value = 2885.1035
$ws.Range("I74").Value = 957.75
$ws.Range("J74").Value = 7168.1113
$ws.Range("K74").Value = 957.75
$ws.Range("L74").Value = 7168.1113
$ws.Range("M74").Value = -83.75
$ws.Range("N74").Value = -8916.1113

$ws.Range("H77").Value = 2885.1035
$ws.Range("I77").Value = 957.75
$ws.Range("J77").Value = 7168.1113
$ws.Range("K77").Value = 4788.75
$ws.Range("L77").Value = 35840.5565
$ws.Range("M77").Value = -420.75
$ws.Range("N77").Value = -44576.5565

$ws.Range("H110").Value = 630.34784
$ws.Range("I110").Value = 613.5454999999999
$ws.Range("K110").Value = 613.5454999999999
$ws.Range("M110").Value = 1431.4545

$ws.Range("H116").Value = 5252.676
$ws.Range("I116").Value = 970.2222
$ws.Range("K116").Value = 970.2222
$ws.Range("M116").Value = 1323.7778

$ws.Range("H122").Value = 2913.2778
$ws.Range("I122").Value = 2661.926
$ws.Range("K122").Value = 7985.778
$ws.Range("M122").Value = -5535.778

$ws.Range("H132").Value = 849780.75
$ws.Range("I132").Value = 1319329.1
$ws.Range("J132").Value = 4593.6
$ws.Range("K132").Value = 3957987.3
$ws.Range("L132").Value = 13780.8
$ws.Range("M132").Value = -3955457.3
$ws.Range("N132").Value = -18840.8

$ws.Range("H136").Value = 878842.5
$ws.Range("I136").Value = 927417.1
$ws.Range("K136").Value = 2782251.3
$ws.Range("M136").Value = -2779701.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5252.676
$ws.Range("I3").Value = 970.2222
$ws.Range("K3").Value = 970.2222
$ws.Range("M3").Value = -856.2222

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H82").Value = 50942
$ws.Range("I82").Value = 16660.666
$ws.Range("J82").Value = 80326
$ws.Range("K82").Value = 16660.666
$ws.Range("L82").Value = 80326
$ws.Range("M82").Value = -16277.666
$ws.Range("N82").Value = -81092

$ws.Range("H85").Value = 50942
$ws.Range("I85").Value = 16660.666
$ws.Range("J85").Value = 80326
$ws.Range("K85").Value = 16660.666
$ws.Range("L85").Value = 80326
$ws.Range("M85").Value = -15334.666
$ws.Range("N85").Value = -82978

$ws.Range("H107").Value = 8335524
$ws.Range("I107").Value = 2131.5483
$ws.Range("J107").Value = 37039430
$ws.Range("K107").Value = 2131.5483
$ws.Range("L107").Value = 37039430
$ws.Range("M107").Value = -211.5482999999999
$ws.Range("N107").Value = -37043270

$ws.Range("H134").Value = 800659.25
$ws.Range("I134").Value = 930818.5600000001
$ws.Range("J134").Value = 462245
$ws.Range("K134").Value = 2792455.68
$ws.Range("L134").Value = 1386735
$ws.Range("M134").Value = -2789920.68
$ws.Range("N134").Value = -1391805

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 93570.05499999999
$ws.Range("I31").Value = 126238.92
$ws.Range("J31").Value = 30745.309
$ws.Range("K31").Value = 126238.92
$ws.Range("L31").Value = 30745.309
$ws.Range("M31").Value = -125943.92
$ws.Range("N31").Value = -31335.309

$ws.Range("H34").Value = 93570.05499999999
$ws.Range("I34").Value = 126238.92
$ws.Range("J34").Value = 30745.309
$ws.Range("K34").Value = 126238.92
$ws.Range("L34").Value = 30745.309
$ws.Range("M34").Value = -126036.92
$ws.Range("N34").Value = -31149.309

$ws.Range("H86").Value = 6565.25
$ws.Range("J86").Value = 7469.4287
$ws.Range("L86").Value = 7469.4287
$ws.Range("N86").Value = -9715.4287

$ws.Range("H89").Value = 6565.25
$ws.Range("J89").Value = 7469.4287
$ws.Range("L89").Value = 37347.14350000001
$ws.Range("N89").Value = -48579.14350000001

$ws.Range("H132").Value = 18967216
$ws.Range("I132").Value = 21279068
$ws.Range("J132").Value = 857718.3
$ws.Range("K132").Value = 63837204
$ws.Range("L132").Value = 2573154.9
$ws.Range("M132").Value = -63834674
$ws.Range("N132").Value = -2578214.9

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 90887
$ws.Range("J37").Value = 90887
$ws.Range("L37").Value = 272661
$ws.Range("N37").Value = -272885

$ws.Range("H87").Value = 14850.23
$ws.Range("I87").Value = 10332.111
$ws.Range("K87").Value = 30996.333
$ws.Range("M87").Value = -29748.333

$ws.Range("H90").Value = 14850.23
$ws.Range("I90").Value = 10332.111
$ws.Range("K90").Value = 92988.99900000001
$ws.Range("M90").Value = -86748.99900000001

$ws.Range("H129").Value = 1159.8572
$ws.Range("I129").Value = 623.9
$ws.Range("J129").Value = 2499.75
$ws.Range("K129").Value = 1871.7
$ws.Range("L129").Value = 7499.25
$ws.Range("M129").Value = 3128.3
$ws.Range("N129").Value = -17499.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3481.625
$ws.Range("I102").Value = 2907.6287
$ws.Range("K102").Value = 2907.6287
$ws.Range("M102").Value = -1285.6287

$ws.Range("H113").Value = 2399.5
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H122").Value = 30926.236
$ws.Range("I122").Value = 44090.043
$ws.Range("J122").Value = 8359.714
$ws.Range("K122").Value = 132270.129
$ws.Range("L122").Value = 25079.142
$ws.Range("M122").Value = -129820.129
$ws.Range("N122").Value = -29979.142

$ws.Range("H126").Value = 836152.6
$ws.Range("I126").Value = 1391186.5
$ws.Range("J126").Value = 3601.75
$ws.Range("K126").Value = 4173559.5
$ws.Range("L126").Value = 10805.25
$ws.Range("M126").Value = -4171089.5
$ws.Range("N126").Value = -15745.25

$ws.Range("H132").Value = 33743824
$ws.Range("I132").Value = 42177148
$ws.Range("J132").Value = 10532.833
$ws.Range("K132").Value = 126531444
$ws.Range("L132").Value = 31598.499
$ws.Range("M132").Value = -126528914
$ws.Range("N132").Value = -36658.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3827.8845
$ws.Range("I7").Value = 3615.476
$ws.Range("K7").Value = 3615.476
$ws.Range("M7").Value = -3503.476

$ws.Range("H61").Value = 11112175
$ws.Range("I61").Value = 1025.125
$ws.Range("J61").Value = 100001380
$ws.Range("K61").Value = 1025.125
$ws.Range("L61").Value = 100001380
$ws.Range("M61").Value = -823.125
$ws.Range("N61").Value = -100001784

$ws.Range("H113").Value = 11112175
$ws.Range("I113").Value = 1025.125
$ws.Range("J113").Value = 100001380
$ws.Range("K113").Value = 1025.125
$ws.Range("L113").Value = 100001380
$ws.Range("M113").Value = 1144.875
$ws.Range("N113").Value = -100005720

$ws.Range("H126").Value = 3827.8845
$ws.Range("I126").Value = 3615.476
$ws.Range("K126").Value = 10846.428
$ws.Range("M126").Value = -8376.428

$ws.Range("H132").Value = 699532.25
$ws.Range("I132").Value = 1089032
$ws.Range("J132").Value = 7088.1113
$ws.Range("K132").Value = 3267096
$ws.Range("L132").Value = 21264.3339
$ws.Range("M132").Value = -3264566
$ws.Range("N132").Value = -26324.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 18746.541
$ws.Range("I100").Value = 1535.4667
$ws.Range("J100").Value = 47431.668
$ws.Range("K100").Value = 3070.9334
$ws.Range("L100").Value = 94863.336
$ws.Range("M100").Value = -2529.9334
$ws.Range("N100").Value = -95945.336

$ws.Range("H113").Value = 1611.4166
$ws.Range("I113").Value = 344.25
$ws.Range("J113").Value = 7947.25
$ws.Range("K113").Value = 1032.75
$ws.Range("L113").Value = 23841.75
$ws.Range("M113").Value = 1137.25
$ws.Range("N113").Value = -28181.75

$ws.Range("H122").Value = 2321.152
$ws.Range("I122").Value = 1978.8857
$ws.Range("K122").Value = 5936.6571
$ws.Range("M122").Value = -3486.6571

$ws.Range("H126").Value = 736.125
$ws.Range("I126").Value = 736.125
$ws.Range("K126").Value = 2208.375
$ws.Range("M126").Value = 261.625

